$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.843.45'
$ws.Range('E2').Value = '  -0.19%  '
$ws.Range('D3').Value = '2.973.73'
$ws.Range('E3').Value = '  -1.18%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '497.42'
$ws.Range('E5').Value = '  -3.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.06'
$ws.Range('E6').Value = '  -2.01%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.427'
$ws.Range('E8').Value = '  -2.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.36'
$ws.Range('E9').Value = '  -1.62%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.107'
$ws.Range('E10').Value = '  -1.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.356'
$ws.Range('E11').Value = '  -0.76%  '
$ws.Range('D12').Value = '3.499.16'
$ws.Range('E12').Value = '  -0.58%  '
$ws.Range('E13').Value = '  -1.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.71'
$ws.Range('E14').Value = '  -1.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000157'
$ws.Range('E15').Value = '  +0.24%  '
$ws.Range('D16').Value = '56.916.68'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('E17').Value = '  +1.84%  '
$ws.Range('D18').Value = '2.972.57'
$ws.Range('E18').Value = '  -0.96%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.59'
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.77'
$ws.Range('E20').Value = '  -1.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '319.61'
$ws.Range('E21').Value = '  -2.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('E23').Value = '  -0.87%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.486'
$ws.Range('E24').Value = '  -0.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.50'
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.161'
$ws.Range('E27').Value = '  -5.91%  '
$ws.Range('D28').Value = '0.0₃0888'
$ws.Range('E28').Value = '  -3.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.57'
$ws.Range('E29').Value = '  -1.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.11'
$ws.Range('E30').Value = '  +0.38%  '
$ws.Range('E31').Value = '  -2.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.15'
$ws.Range('E32').Value = '  -6.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.12'
$ws.Range('E33').Value = '  -2.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '152.45'
$ws.Range('E34').Value = '  -1.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.61'
$ws.Range('E35').Value = '  +0.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.74'
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('E37').Value = '  -2.72%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '23.99'
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0665'
$ws.Range('E39').Value = '  -2.25%  '
$ws.Range('D40').Value = '3.004.33'
$ws.Range('E40').Value = '  -1.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '37.53'
$ws.Range('E41').Value = '  +1.02%  '
$ws.Range('E42').Value = '  +0.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.73'
$ws.Range('E43').Value = '  +0.79%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.639'
$ws.Range('E44').Value = '  -1.39%  '
$ws.Range('D45').Value = '2.193.74'
$ws.Range('E45').Value = '  -4.54%  '
$ws.Range('E46').Value = '  -3.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.945'
$ws.Range('E47').Value = '  -6.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.92'
$ws.Range('E48').Value = '  +0.34%  '
$ws.Range('E49').Value = '  -1.98%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.09'
$ws.Range('E50').Value = '  -2.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.79'
$ws.Range('E51').Value = '  -9.04%  '
